$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03-t03-activity-deleted")

$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "Finish"
